# Update "想去人数" (want-to-go count) figures in the F column
# for the 展览 (Exhibitions) sheet and the 全部类型 (All types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 147
$ws1.Range("F3").Value = 23
$ws1.Range("F4").Value = 224
$ws1.Range("F5").Value = 3788
$ws1.Range("F6").Value = 390
$ws1.Range("F7").Value = 23

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 147
$ws4.Range("F3").Value = 23
$ws4.Range("F4").Value = 224
$ws4.Range("F5").Value = 3788
$ws4.Range("F6").Value = 390
$ws4.Range("F9").Value = 23
